# The document has three inline logo pictures living in the headers/footers
# (two copies of the Pearson logo in the footers, one BTEC logo in the
# header). This rename swaps their internal drawing "Name" (the
# wp:docPr / pic:cNvPr "name" attribute, i.e. the picture's Name as shown
# in Word's Selection Pane) without touching anything else about the
# picture (size, position, alt text, embedded media, etc.):
#
#   footer Pearson logos : image2.png -> image1.png
#   header BTEC logo      : image1.jpg -> image2.jpg
#
# InlineShape does not expose a settable Name in the Word object model,
# so each picture is briefly converted to a floating Shape (which does
# expose .Name), renamed, then converted back to an inline shape so the
# layout/anchoring in the document is left exactly as it was.

function Rename-InlineShape($inlineShape, $newName) {
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape()
}

$d = $word.ActiveDocument
$sec = $d.Sections.First

# Footers 1 and 2 each hold one Pearson logo inline picture.
for ($i = 1; $i -le 3; $i++) {
    $f = $sec.Footers.Item($i)
    if ($f.Exists) {
        for ($j = 1; $j -le $f.Range.InlineShapes.Count; $j++) {
            $shp = $f.Range.InlineShapes.Item($j)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                Rename-InlineShape $shp "image1.png"
            }
        }
    }
}

# One of the headers holds the BTEC logo inline picture.
for ($i = 1; $i -le 3; $i++) {
    $h = $sec.Headers.Item($i)
    if ($h.Exists) {
        for ($j = 1; $j -le $h.Range.InlineShapes.Count; $j++) {
            $shp = $h.Range.InlineShapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                Rename-InlineShape $shp "image2.jpg"
            }
        }
    }
}
